$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three distinct data rows (well "JY41-5", facies 3) that get
# appended three times in a row (398-400, 401-403, 404-406).
$blockRows = @(
    @(2615.625, 3, 213.388, 226.21,  2.613, 7.306, 40.447, 42.848, 2.431),
    @(2615.75,  3, 203.697, 228.734, 2.61,  7.849, 30.651, 32.818, 2.231),
    @(2615.875, 3, 210.302, 234.457, 2.569, 9.559, 21.823, 23.516, 2.081)
)

$startRow = 398
for ($block = 0; $block -lt 3; $block++) {
    for ($i = 0; $i -lt 3; $i++) {
        $r = $startRow + ($block * 3) + $i
        $data = $blockRows[$i]

        $ws.Cells.Item($r, 1).Value = "JY41-5"
        $ws.Cells.Item($r, 2).Value = $data[0]
        $ws.Cells.Item($r, 3).Value = $data[1]
        $ws.Cells.Item($r, 4).Value = $data[2]
        $ws.Cells.Item($r, 5).Value = $data[3]
        $ws.Cells.Item($r, 6).Value = $data[4]
        $ws.Cells.Item($r, 7).Value = $data[5]
        $ws.Cells.Item($r, 8).Value = $data[6]
        $ws.Cells.Item($r, 9).Value = $data[7]
        $ws.Cells.Item($r, 10).Value = $data[8]
    }
}

# Re-apply the autofilter so its range grows to cover the newly added
# rows (C1:C397 -> C1:C406).
$ws.AutoFilterMode = $false
$ws.Range("C1:C406").AutoFilter()

# The hidden workbook-level _FilterDatabase name also tracks the
# autofilter range; make sure it is widened to match.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$C`$1:`$C`$406"
    }
}

# Update the active selection to mirror the post-edit cursor position.
$ws.Range("D412").Select()
